# Auto-generated script to reorder email lists in column G
# (same set of emails per row, just in a different order - matches the commit diff)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$ws.Range("G2").Value = 'mennatulla.medhat@med.asu.edu.eg, servinaz@med.asu.edu.eg'
$ws.Range("G3").Value = 'shaimaa.ahmed@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg, servinaz@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg'
$ws.Range("G4").Value = 'Veronia.rafat@med.asu.edu.eg, servinaz@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg, hananragab@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg, nesmadrahim@med.asu.edu.eg'
$ws.Range("G5").Value = 'heba@med.asu.edu.eg, gehanadel@med.asu.edu.eg, System, rana.abozaid@med.asu.edu.eg, servinaz@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg'
$ws.Range("G6").Value = 'Amira.Sobhy@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg'
$ws.Range("G8").Value = 'NadaMohamed@med.asu.edu.eg, AbeerRagheb@med.asu.edu.eg, menna-alah.mohamed@asu.edu.eg'
$ws.Range("G9").Value = 'wessam.atef@med.asu.edu.eg, norhan.mohamed@med.asu.edu.eg, Safa.hany@med.asu.edu.eg, Sara_nabil@med.asu.edu.eg'
$ws.Range("G10").Value = 'aml.awwad@med.asu.edu.eg, aya.saeed@med.asu.edu.eg, Omnia.Mohammed@med.asu.edu.eg, wessam.atef@med.asu.edu.eg, Sara_nabil@med.asu.edu.eg'
$ws.Range("G12").Value = 'Eman.m.abosakaya@med.asu.edu.eg, basma.hamed@med.asu.edu.eg, esraa.mostafa@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, marwa_mustafa@med.asu.edu.eg, merna.said@med.asu.edu.eg, mernamahrous@med.asu.edu.eg, amany.raafat@med.asu.edu.eg, maryam.ahmed@med.asu.edu.eg'
$ws.Range("G13").Value = 'maimustafa@med.asu.edu.eg, dina.adel@med.asu.edu.eg, nourhan.osama@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, Eman.m.elsaid@med.asu.edu.eg, amira.m.ibrahim@med.asu.edu.eg'
$ws.Range("G15").Value = 'mariam.youssif.std@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg, nourhanmamdouh@med.asu.edu.eg, Mayarembaby@med.asu.edu.eg'
$ws.Range("G16").Value = 'mohamed.saleem@med.asu.edu.eg, Alshimaa_khaled@med.asu.edu.eg, Mayarembaby@med.asu.edu.eg'
$ws.Range("G25").Value = 'youstina.gamil@med.asu.edu.eg, mona.I.hussein@med.asu.edu.eg'
$ws.Range("G27").Value = 'nourhanmohamed@med.asu.edu.eg, marian.samir@med.asu.edu.eg, randa.rabea@med.asu.edu.eg, System'
$ws.Range("G29").Value = 'aya.hanafy@med.asu.edu.eg, eman.samir@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, nahla91@med.asu.edu.eg, marinasorial@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg, marina_atef@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg'
$ws.Range("G30").Value = 'nahla91@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg'
$ws.Range("G33").Value = 'mennatulla.medhat@med.asu.edu.eg, servinaz@med.asu.edu.eg'
$ws.Range("G34").Value = 'shaimaa.ahmed@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg, servinaz@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg'
$ws.Range("G35").Value = 'gehanadel@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg'
$ws.Range("G36").Value = 'heba@med.asu.edu.eg, gehanadel@med.asu.edu.eg, System, rana.abozaid@med.asu.edu.eg, servinaz@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg'
$ws.Range("G37").Value = 'Amira.Sobhy@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg'
$ws.Range("G38").Value = 'AbeerRagheb@med.asu.edu.eg, Fatmaelhady@med.asu.edu.eg'
$ws.Range("G39").Value = 'NadaMohamed@med.asu.edu.eg, menna-alah.mohamed@asu.edu.eg'
$ws.Range("G40").Value = 'wessam.atef@med.asu.edu.eg, norhan.mohamed@med.asu.edu.eg, Safa.hany@med.asu.edu.eg, Sara_nabil@med.asu.edu.eg'
$ws.Range("G41").Value = 'aml.awwad@med.asu.edu.eg, norhan.mohamed@med.asu.edu.eg, Sara_nabil@med.asu.edu.eg'
$ws.Range("G43").Value = 'Eman.m.abosakaya@med.asu.edu.eg, basma.hamed@med.asu.edu.eg, esraa.mostafa@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, marwa_mustafa@med.asu.edu.eg, merna.said@med.asu.edu.eg, mernamahrous@med.asu.edu.eg, amany.raafat@med.asu.edu.eg, maryam.ahmed@med.asu.edu.eg'
$ws.Range("G44").Value = 'maimustafa@med.asu.edu.eg, dina.adel@med.asu.edu.eg, nourhan.osama@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, Eman.m.elsaid@med.asu.edu.eg, amira.m.ibrahim@med.asu.edu.eg'
$ws.Range("G46").Value = 'mohamed.saleem@med.asu.edu.eg, Alshimaa_khaled@med.asu.edu.eg, afnan.fares@med.asu.edu.eg'
$ws.Range("G56").Value = 'youstina.gamil@med.asu.edu.eg, mona.I.hussein@med.asu.edu.eg'
$ws.Range("G58").Value = 'nourhanmohamed@med.asu.edu.eg, marian.samir@med.asu.edu.eg, randa.rabea@med.asu.edu.eg, System'
$ws.Range("G60").Value = 'aya.hanafy@med.asu.edu.eg, eman.samir@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, nahla91@med.asu.edu.eg, marinasorial@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg, marina_atef@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg'
$ws.Range("G61").Value = 'nahla91@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg'
$ws.Range("G64").Value = 'Veronia.rafat@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg'
$ws.Range("G65").Value = 'Amira.Sobhy@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg'
$ws.Range("G66").Value = 'Veronia.rafat@med.asu.edu.eg, servinaz@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg, hananragab@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg, nesmadrahim@med.asu.edu.eg'
$ws.Range("G67").Value = 'gehanadel@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, servinaz@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg'
$ws.Range("G68").Value = 'majorelle.magdy@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg, alshimaa.atef@med.asu.edu.egm, rana.abozaid@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg'
$ws.Range("G69").Value = 'Fatmaelhady@med.asu.edu.eg, Kerelos.zareef@med.asu.edu.eg'
$ws.Range("G70").Value = 'NadaMohamed@med.asu.edu.eg, AbeerRagheb@med.asu.edu.eg, menna-alah.mohamed@asu.edu.eg'
$ws.Range("G71").Value = 'Omnia.Mohammed@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg, Safa.hany@med.asu.edu.eg, Sara_nabil@med.asu.edu.eg'
$ws.Range("G72").Value = 'aml.awwad@med.asu.edu.eg, aya.saeed@med.asu.edu.eg, Omnia.Mohammed@med.asu.edu.eg, wessam.atef@med.asu.edu.eg, Sara_nabil@med.asu.edu.eg'
$ws.Range("G74").Value = 'maimustafa@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, mernamahrous@med.asu.edu.eg, merna.said@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg, amany.raafat@med.asu.edu.eg'
$ws.Range("G75").Value = 'maryam.ahmed@med.asu.edu.eg, amany.raafat@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg'
$ws.Range("G76").Value = 'mariam.youssif.std@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg, nourhanmamdouh@med.asu.edu.eg, Mayarembaby@med.asu.edu.eg'
$ws.Range("G77").Value = 'mohamed.saleem@med.asu.edu.eg, Alshimaa_khaled@med.asu.edu.eg, afnan.fares@med.asu.edu.eg'
$ws.Range("G91").Value = 'aya.hanafy@med.asu.edu.eg, eman.samir@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, nahla91@med.asu.edu.eg, marinasorial@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg, marina_atef@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg'
$ws.Range("G92").Value = 'nahla91@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg'
$ws.Range("G95").Value = 'Veronia.rafat@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, servinaz@med.asu.edu.eg'
$ws.Range("G96").Value = 'nourhan.mahmoud@med.asu.edu.eg, servinaz@med.asu.edu.eg, alshimaa.atef@med.asu.edu.egm'
$ws.Range("G97").Value = 'gehanadel@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, servinaz@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg'
$ws.Range("G98").Value = 'heba@med.asu.edu.eg, gehanadel@med.asu.edu.eg, nahla.nagiub@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg'
$ws.Range("G101").Value = 'NadaMohamed@med.asu.edu.eg, AbeerRagheb@med.asu.edu.eg, menna-alah.mohamed@asu.edu.eg'
$ws.Range("G102").Value = 'aml.awwad@med.asu.edu.eg, norhan.mohamed@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg, Safa.hany@med.asu.edu.eg'
$ws.Range("G103").Value = 'aml.awwad@med.asu.edu.eg, norhan.mohamed@med.asu.edu.eg, Sara_nabil@med.asu.edu.eg'
$ws.Range("G105").Value = 'Eman.m.abosakaya@med.asu.edu.eg, basma.hamed@med.asu.edu.eg, esraa.mostafa@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, marwa_mustafa@med.asu.edu.eg, merna.said@med.asu.edu.eg, mernamahrous@med.asu.edu.eg, amany.raafat@med.asu.edu.eg, maryam.ahmed@med.asu.edu.eg'
$ws.Range("G106").Value = 'dina.adel@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, nadia.mostafa@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg'
$ws.Range("G108").Value = 'mariam.youssif.std@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg, nourhanmamdouh@med.asu.edu.eg, Mayarembaby@med.asu.edu.eg'
$ws.Range("G109").Value = 'mohamed.saleem@med.asu.edu.eg, Alshimaa_khaled@med.asu.edu.eg, Mayarembaby@med.asu.edu.eg'
$ws.Range("G120").Value = 'nourhanmohamed@med.asu.edu.eg, marian.samir@med.asu.edu.eg, dr_ahmed_mostafa@med.asu.edu.eg, nourhan.mostafa@med.asu.edu.eg, afaf.abdallah@med.asu.edu.eg'
$ws.Range("G122").Value = 'eman.samir@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg'
$ws.Range("G123").Value = 'nahla91@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg'
$ws.Range("G126").Value = 'Veronia.rafat@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, servinaz@med.asu.edu.eg'
$ws.Range("G127").Value = 'nourhan.mahmoud@med.asu.edu.eg, servinaz@med.asu.edu.eg, alshimaa.atef@med.asu.edu.egm'
$ws.Range("G128").Value = 'gehanadel@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, servinaz@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg'
$ws.Range("G129").Value = 'heba@med.asu.edu.eg, gehanadel@med.asu.edu.eg, nahla.nagiub@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg'
$ws.Range("G131").Value = 'Fatmaelhady@med.asu.edu.eg, Kerelos.zareef@med.asu.edu.eg'
$ws.Range("G132").Value = 'NadaMohamed@med.asu.edu.eg, AbeerRagheb@med.asu.edu.eg, menna-alah.mohamed@asu.edu.eg'
$ws.Range("G133").Value = 'aml.awwad@med.asu.edu.eg, norhan.mohamed@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg, Safa.hany@med.asu.edu.eg'
$ws.Range("G134").Value = 'aml.awwad@med.asu.edu.eg, aya.saeed@med.asu.edu.eg'
$ws.Range("G136").Value = 'Eman.m.abosakaya@med.asu.edu.eg, basma.hamed@med.asu.edu.eg, esraa.mostafa@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, marwa_mustafa@med.asu.edu.eg, merna.said@med.asu.edu.eg, mernamahrous@med.asu.edu.eg, amany.raafat@med.asu.edu.eg, maryam.ahmed@med.asu.edu.eg'
$ws.Range("G137").Value = 'maimustafa@med.asu.edu.eg, dina.adel@med.asu.edu.eg, nourhan.osama@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, Eman.m.elsaid@med.asu.edu.eg, amira.m.ibrahim@med.asu.edu.eg'
$ws.Range("G138").Value = 'mariam.youssif.std@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg, nourhanmamdouh@med.asu.edu.eg, Mayarembaby@med.asu.edu.eg'
$ws.Range("G139").Value = 'mariam.youssif.std@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg, Mayarembaby@med.asu.edu.eg, Rania.a.youssef@med.asu.edu.eg'
$ws.Range("G140").Value = 'mohamed.saleem@med.asu.edu.eg, Alshimaa_khaled@med.asu.edu.eg, afnan.fares@med.asu.edu.eg'
$ws.Range("G151").Value = 'enas.omran@med.asu.edu.eg, user@user.com, takyas@med.asu.edu.eg, Amr-Saeed@med.asu.edu.eg'
$ws.Range("G153").Value = 'eman.samir@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg'
$ws.Range("G154").Value = 'youstina.magdy@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg, marina_atef@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, eman.samir@med.asu.edu.eg'
$ws.Range("G157").Value = 'majorelle.magdy@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, servinaz@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg'
$ws.Range("G158").Value = 'Amira.Sobhy@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg'
$ws.Range("G159").Value = 'gehanadel@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg'
$ws.Range("G160").Value = 'Amira.Sobhy@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, servinaz@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg'
$ws.Range("G161").Value = 'majorelle.magdy@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg, alshimaa.atef@med.asu.edu.egm, rana.abozaid@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg'
$ws.Range("G162").Value = 'Fatmaelhady@med.asu.edu.eg, Kerelos.zareef@med.asu.edu.eg'
$ws.Range("G164").Value = 'Omnia.Mohammed@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg, Safa.hany@med.asu.edu.eg, Sara_nabil@med.asu.edu.eg'
$ws.Range("G165").Value = 'aml.awwad@med.asu.edu.eg, aya.saeed@med.asu.edu.eg'
$ws.Range("G167").Value = 'maimustafa@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, mernamahrous@med.asu.edu.eg, merna.said@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg, amany.raafat@med.asu.edu.eg'
$ws.Range("G168").Value = 'maryam.ahmed@med.asu.edu.eg, amany.raafat@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg'
$ws.Range("G170").Value = 'mariam.youssif.std@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg, Mayarembaby@med.asu.edu.eg, Rania.a.youssef@med.asu.edu.eg'
$ws.Range("G182").Value = 'nourhanmohamed@med.asu.edu.eg, marian.samir@med.asu.edu.eg, dr_ahmed_mostafa@med.asu.edu.eg, nourhan.mostafa@med.asu.edu.eg, afaf.abdallah@med.asu.edu.eg'
$ws.Range("G184").Value = 'yasmin.m.senosy@med.asu.edu.eg, nahla91@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg'
$ws.Range("G185").Value = 'youstina.magdy@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg, marina_atef@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, eman.samir@med.asu.edu.eg'
